$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timestamp value applied to all data rows (2-9) after re-run
$newTimestamp = 45841.64214048738

# Updated category strings (B column) reflecting stricter Coding threshold
# (>=5 prompts & cos>=0.45) applied in Q06, which also triggered a re-run of Q07-Q09
$ws.Range("B2").Value = "Writing & professional communication|Coding - programming help|Study revision - exam prep|Other"
$ws.Range("B4").Value = "Other"
$ws.Range("B5").Value = "Writing & professional communication|Coding - programming help|Other"
$ws.Range("B7").Value = "Writing & professional communication|Brainstorming & personal ideas - fun|Language practice or translation|Other"
$ws.Range("B8").Value = "Writing & professional communication|Coding - programming help|Other"
$ws.Range("B9").Value = "Writing & professional communication|Brainstorming & personal ideas - fun|Coding - programming help|Language practice or translation|Other"

# Update timestamps for all data rows (2-9)
$ws.Range("D2").Value = $newTimestamp
$ws.Range("D3").Value = $newTimestamp
$ws.Range("D4").Value = $newTimestamp
$ws.Range("D5").Value = $newTimestamp
$ws.Range("D6").Value = $newTimestamp
$ws.Range("D7").Value = $newTimestamp
$ws.Range("D8").Value = $newTimestamp
$ws.Range("D9").Value = $newTimestamp

$wb.Save()
